$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# Step 1: Insert the two new header rows (this shifts cell values, styles,
# merged cells and data validations down automatically).
$ws.Rows.Item(89).Insert()
$ws.Rows.Item(99).Insert()

# Step 2: Fix up the two newly inserted rows: copy header formatting from an
# existing section-header row (row 85), set their text, and drop the spurious
# B-column cell that Insert() created (header rows only use column A).
$ws.Range("A85").Copy()
$ws.Range("A89").PasteSpecial(-4122)
$ws.Range("A89").Value2 = "KVM RAM"
$ws.Range("B89").Clear()

$ws.Range("A85").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$ws.Range("A99").Value2 = "KVM CPU"
$ws.Range("B99").Clear()
$ws.Application.CutCopyMode = 0

# Step 3: Update the label text for every row per the new titles.
$ws.Range("A84").Value2 = 'NSGV Disk Size'
$ws.Range("A85").Value2 = 'KVM and VCenter RAM'
$ws.Range("A86").Value2 = 'KVM VSD RAM'
$ws.Range("A87").Value2 = 'KVM VSC RAM'
$ws.Range("A88").Value2 = 'KVM VSTAT RAM'
$ws.Range("A89").Value2 = 'KVM RAM'
$ws.Range("A90").Value2 = 'KVM VCIN RAM'
$ws.Range("A91").Value2 = 'KVM NUH RAM'
$ws.Range("A92").Value2 = 'KVM Webfilter RAM'
$ws.Range("A93").Value2 = 'KVM Portal RAM'
$ws.Range("A94").Value2 = 'KVM and VCenter CPU'
$ws.Range("A95").Value2 = 'KVM VSD CPU cores'
$ws.Range("A96").Value2 = 'KVM VSC CPU cores'
$ws.Range("A97").Value2 = 'KVM VSTAT CPU cores'
$ws.Range("A98").Value2 = 'KVM VNSUTIL CPU cores'
$ws.Range("A99").Value2 = 'KVM CPU'
$ws.Range("A100").Value2 = 'KVM NUH CPU cores'
$ws.Range("A101").Value2 = 'KVM VCIN CPU cores'
$ws.Range("A102").Value2 = 'KVM Portal VM CPU cores'
$ws.Range("A103").Value2 = 'KVM Webfilter VM CPU cores'
$ws.Range("A104").Value2 = 'Authentication'
$ws.Range("A105").Value2 = 'VSD Architect URL'
$ws.Range("A106").Value2 = 'VSD Enterprise'
$ws.Range("A107").Value2 = 'VCIN URL'
$ws.Range("A108").Value2 = 'VCIN Enterprise'
$ws.Range("A109").Value2 = 'Hooks'
$ws.Range("A110").Value2 = 'hooks'
$ws.Range("A111").Value2 = 'skip actions'
$ws.Range("A112").Value2 = 'Component Health Report Email Options'
$ws.Range("A113").Value2 = 'Health Report SMTP Server'
$ws.Range("A114").Value2 = 'Health Report SMTP Server Port'
$ws.Range("A115").Value2 = 'Health Report Email From Address'
$ws.Range("A116").Value2 = 'Health Report Destination Email Address(es)'
$ws.Range("A117").Value2 = 'VSD Monit Email Alerts Configuration'
$ws.Range("A118").Value2 = 'VSD Monit Mail Server'
$ws.Range("A119").Value2 = 'VSD Monit Mail Server Port'
$ws.Range("A120").Value2 = 'VSD Monit Mail Server Encryption Type'
$ws.Range("A121").Value2 = 'Use VSD Monit Eventqueue'
$ws.Range("A122").Value2 = 'VSD Monit From Email Address'
$ws.Range("A123").Value2 = 'VSD Monit Reply-To Email Address'
$ws.Range("A124").Value2 = 'VSD Monit Email Alert Subject'
$ws.Range("A125").Value2 = 'VSD Monit Email Alert Message'
$ws.Range("A126").Value2 = 'VSD Monit Destination Email Address'
$ws.Range("A127").Value2 = 'VSD Monit Only Alert On'
$ws.Range("A128").Value2 = 'VSD Monit Do Not Alert On'
$ws.Range("A129").Value2 = 'VSD In-place upgrade during Install '

# Step 4: Comments did not shift automatically with the row inserts, so fix
# them up explicitly: remove stale comments left on header rows, update the
# text on comments that stayed in place, and add comments that are now on
# rows that did not have one before.
# 4a. delete stale comments on rows that must not have a comment any more
$ws.Range("A89").Comment.Delete()
$ws.Range("A94").Comment.Delete()
$ws.Range("A99").Comment.Delete()
$ws.Range("A104").Comment.Delete()
$ws.Range("A109").Comment.Delete()
$ws.Range("A112").Comment.Delete()
$ws.Range("A117").Comment.Delete()

# 4b. update text of comments staying on the same row index
$ws.Range("A84").Comment.Text('Amount of NSGV disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments should not modify this value. [default: 4]')
$ws.Range("A86").Comment.Text('For KVM and VCenter deployments: amount of VSD RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$ws.Range("A87").Comment.Text('For KVM and VCenter deployments: amount of VSC RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 4]')
$ws.Range("A88").Comment.Text('For KVM and VCenter deployments: amount of VSTAT RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 16]')
$ws.Range("A90").Comment.Text('Amount of VCIN RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$ws.Range("A91").Comment.Text('Amount of NUH RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]')
$ws.Range("A92").Comment.Text('Amount of Webfilter RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]')
$ws.Range("A95").Comment.Text('For KVM and VCenter deployments: number of CPU''s for VSD. [default: 6]')
$ws.Range("A96").Comment.Text('For KVM and VCenter deployments: number of CPU''s for VSC. [default: 6]')
$ws.Range("A97").Comment.Text('For KVM and VCenter deployments: number of CPU''s for VSTAT. [default: 6]')
$ws.Range("A98").Comment.Text('For KVM and VCenter deployments: number of CPU''s for VNSUTIL. [default: 2]')
$ws.Range("A100").Comment.Text('Number of CPU''s for NUH. Valid only for KVM deployments [default: 2]')
$ws.Range("A101").Comment.Text('Number of CPU''s for VCIN. Valid only for KVM deployments [default: 6]')
$ws.Range("A103").Comment.Text('Number of CPU''s for Webfilter vm. Valid only for KVM deployments [default: 2]')
$ws.Range("A105").Comment.Text('VSD Architect URL. Required for tasks during Upgrade, Health Checks etc [default: https://(vsd_fqdn):8443]')
$ws.Range("A106").Comment.Text('Enterprise name used for authentication with VSD Architect. Required for tasks during Upgrade, Health Checks etc [default: csp]')
$ws.Range("A108").Comment.Text('Enterprise name used for authentication with VCIN. Required for tasks like VRS-E upgrade (through VCIN) [default: csp]')
$ws.Range("A111").Comment.Text('Skip tasks and playbooks (List items separated by comma.)')
$ws.Range("A113").Comment.Text('Address of SMTP server to be used if emailing health results')
$ws.Range("A114").Comment.Text('Port to be used on the SMTP Server [default: 25]')
$ws.Range("A116").Comment.Text('List of destination email addresses (List items separated by comma.)')
$ws.Range("A118").Comment.Text('Address of the mail server to be used to receive monit alerts via email')
$ws.Range("A119").Comment.Text('Port on mail server to be used for monit alerts [default: 25]')
$ws.Range("A120").Comment.Text('Encryption to be used when sending monit alerts via email')
$ws.Range("A121").Comment.Text('Enables use of monit eventqueue to store alerts if email alerts fail to send [default: True]')
$ws.Range("A122").Comment.Text('Email address from which monit alerts will be sent')
$ws.Range("A123").Comment.Text('Email address to reply to monit alert emails')
$ws.Range("A124").Comment.Text('Email subject for alert emails. Overrides monit default alert subject')
$ws.Range("A125").Comment.Text('Email message for alert emails. Overrides monit default alert message')
$ws.Range("A126").Comment.Text('Destination email address for monit alerts')
$ws.Range("A127").Comment.Text('Specific events for which alerts should be sent. One string can be used to hold multiple events, separated by commas')

# 4c. add comments on rows that do not have one yet
$ws.Range("A93").AddComment('Amount of Portal RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$ws.Range("A102").AddComment('Number of CPU''s for Portal vm. Valid only for KVM deployments [default: 6]')
$ws.Range("A107").AddComment('VCIN URL used for API interaction. Required for tasks like VRS-E upgrade (through VCIN) [default: https://(vcin_ip_address):8443]')
$ws.Range("A110").AddComment('List of hooks files (List items separated by comma.)')
$ws.Range("A115").AddComment('Email address from which health report will be sent')
$ws.Range("A128").AddComment('Events for which alerts should not be sent. One string can be used to hold multiple events, separated by commas')
$ws.Range("A129").AddComment('Allowing VSD in-place upgrade during Installation [default: False]')
